$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update sheet (tab) name and title to reflect new "through" date
$ws.Name = "Through 2021-12-10"

# Update the December row label shared string
$ws.Range("A13").Value = "December (through 12-10)"

# Update the 2021 column (H) values affected by the new data
$ws.Range("H11").Value = 196
$ws.Range("H13").Value = 75
$ws.Range("H14").Value = 1719

# Update row 13 (December) values for years 2015-2020 (B-G)
$ws.Range("B13").Value = 9
$ws.Range("C13").Value = 28
$ws.Range("E13").Value = 24
$ws.Range("F13").Value = 14
$ws.Range("G13").Value = 48

# Update row 14 (Total) values for years 2015-2020 (B-G)
$ws.Range("B14").Value = 300
$ws.Range("C14").Value = 591
$ws.Range("E14").Value = 706
$ws.Range("F14").Value = 548
$ws.Range("G14").Value = 1312
